$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: reorder/rename the lecture entries
$ws.Range("B2").Value = "CD222-sec-Hall 1`nCD222-lab-Hall 2`nCD222-Dr. Wendy Stokes-Hall 5"

# Update C2: the "sec" entry becomes the Carlos Langworth entry
$ws.Range("C2").Value = "CD222-Carlos Langworth-Hall 1"
